$wb = $excel.ActiveWorkbook

# Duplicate the "September 2022" sheet, placing the copy before it.
# This becomes the new "October 2022" sheet, inheriting headers/formatting.
$sept = $wb.Worksheets.Item("September 2022")
$sept.Copy($sept)

# Re-fetch fresh references by name (the sheet collection shifted).
$newSheet = $wb.Worksheets.Item("September 2022 (2)")
$newSheet.Name = "October 2022"
$sept = $wb.Worksheets.Item("September 2022")

# Update the flagged-channel entries for October.
$newSheet.Range("B2").Value = "LBC_m10_c37_highgain"
$newSheet.Range("B2").Font.Color = 10498160

$newSheet.Range("B3").Value = "LBC_m20_c37_highgain"
$newSheet.Range("B4").Value = "LBC_m20_c37_lowgain"
$newSheet.Range("B5").ClearContents()

# Update selection on the (no longer active) September 2022 sheet.
$sept.Select()
$sept.Range("A1:B6").Select()

# Restore October 2022 as the active/visible tab with its own selection.
$newSheet.Select()
$newSheet.Range("B4").Select()
